$d = $word.ActiveDocument

# Remove the stray " souhvězdí " between "zobrazují" and "Herkules"
# (e.g. "zobrazují souhvězdí Herkules" -> "zobrazujíHerkules")
$d.Content.Find.Execute(
    "zobrazují souhvězdí Herkules",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "zobrazujíHerkules",
    2
)
